# Auto-generated script applying cell-level numeric updates
# described by the OOXML diff for Sheets/Halicarnassus_Profits.xlsx
# (workbook tabs: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 85.25
$ws.Range("I6").Value = 20.5
$ws.Range("J6").Value = 150
$ws.Range("K6").Value = 61.5
$ws.Range("L6").Value = 450
$ws.Range("M6").Value = 50.5
$ws.Range("N6").Value = -674

$ws.Range("H64").Value = 9424.75
$ws.Range("J64").Value = 9599.75
$ws.Range("L64").Value = 9599.75
$ws.Range("N64").Value = -10095.75

$ws.Range("H67").Value = 9424.75
$ws.Range("J67").Value = 9599.75
$ws.Range("L67").Value = 9599.75
$ws.Range("N67").Value = -11315.75

$ws.Range("H74").Value = 19999.334
$ws.Range("I74").Value = 15998.8
$ws.Range("K74").Value = 15998.8
$ws.Range("M74").Value = -15062.8

$ws.Range("H77").Value = 19999.334
$ws.Range("I77").Value = 15998.8
$ws.Range("K77").Value = 79994
$ws.Range("M77").Value = -75314

$ws.Range("H118").Value = 1624.75
$ws.Range("J118").Value = 5000
$ws.Range("L118").Value = 15000
$ws.Range("N118").Value = -18314

$ws.Range("H121").Value = 573.36365
$ws.Range("J121").Value = 573.36365
$ws.Range("L121").Value = 1720.09095
$ws.Range("N121").Value = -5214.09095

$ws.Range("H138").Value = 4195.9
$ws.Range("J138").Value = 4195.9
$ws.Range("L138").Value = 12587.7
$ws.Range("N138").Value = -22867.7

$ws.Range("H141").Value = 2634.2727
$ws.Range("I141").Value = 1834.625
$ws.Range("K141").Value = 5503.875
$ws.Range("M141").Value = -323.875


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 150.73685
$ws.Range("I22").Value = 150.73685
$ws.Range("K22").Value = 150.73685
$ws.Range("M22").Value = 148.26315

$ws.Range("H61").Value = 5615.5
$ws.Range("I61").Value = 5205.467
$ws.Range("J61").Value = 7665.6665
$ws.Range("K61").Value = 5205.467
$ws.Range("L61").Value = 7665.6665
$ws.Range("M61").Value = -4993.467
$ws.Range("N61").Value = -8089.6665

$ws.Range("H63").Value = 3257.8
$ws.Range("I63").Value = 3257.8
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 3257.8
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -2571.8
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 3257.8
$ws.Range("I66").Value = 3257.8
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 16289
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -12857
$ws.Range("N66").ClearContents()

$ws.Range("H132").Value = 3535.5715
$ws.Range("I132").Value = 2790
$ws.Range("J132").Value = 5399.5
$ws.Range("K132").Value = 8370
$ws.Range("L132").Value = 16198.5
$ws.Range("M132").Value = -5840
$ws.Range("N132").Value = -21258.5

$ws.Range("H136").Value = 5615.5
$ws.Range("I136").Value = 5205.467
$ws.Range("J136").Value = 7665.6665
$ws.Range("K136").Value = 15616.401
$ws.Range("L136").Value = 22996.9995
$ws.Range("M136").Value = -13066.401
$ws.Range("N136").Value = -28096.9995


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 8754
$ws.Range("I20").Value = 2508
$ws.Range("K20").Value = 2508
$ws.Range("M20").Value = -2261

$ws.Range("H22").Value = 354.2353
$ws.Range("J22").Value = 399.25
$ws.Range("L22").Value = 399.25
$ws.Range("N22").Value = -745.25

$ws.Range("H134").Value = 7748.75
$ws.Range("I134").Value = 6998.6665
$ws.Range("K134").Value = 20995.9995
$ws.Range("M134").Value = -18460.9995


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2971.5715
$ws.Range("J16").Value = 3663.3333
$ws.Range("L16").Value = 3663.3333
$ws.Range("N16").Value = -4237.3333

$ws.Range("H31").Value = 9166.147999999999
$ws.Range("J31").Value = 10037.692
$ws.Range("L31").Value = 10037.692
$ws.Range("N31").Value = -10627.692

$ws.Range("H34").Value = 9166.147999999999
$ws.Range("J34").Value = 10037.692
$ws.Range("L34").Value = 10037.692
$ws.Range("N34").Value = -10441.692

$ws.Range("H58").Value = 1961.875
$ws.Range("I58").Value = 1965.5
$ws.Range("J58").Value = 1951
$ws.Range("K58").Value = 1965.5
$ws.Range("L58").Value = 1951
$ws.Range("M58").Value = -1762.5
$ws.Range("N58").Value = -2357

$ws.Range("H113").Value = 2971.5715
$ws.Range("J113").Value = 3663.3333
$ws.Range("L113").Value = 3663.3333
$ws.Range("N113").Value = -8003.3333

$ws.Range("H136").Value = 1961.875
$ws.Range("I136").Value = 1965.5
$ws.Range("J136").Value = 1951
$ws.Range("K136").Value = 5896.5
$ws.Range("L136").Value = 5853
$ws.Range("M136").Value = -3346.5
$ws.Range("N136").Value = -10953


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 47.333332
$ws.Range("I2").Value = 47.5
$ws.Range("K2").Value = 285
$ws.Range("M2").Value = -172

$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9504.75
$ws.Range("I70").Value = 8000
$ws.Range("K70").Value = 8000
$ws.Range("M70").Value = -7730

$ws.Range("H73").Value = 9504.75
$ws.Range("I73").Value = 8000
$ws.Range("K73").Value = 8000
$ws.Range("M73").Value = -7064

$ws.Range("H80").Value = 4499.75
$ws.Range("J80").Value = 4666.3335
$ws.Range("L80").Value = 4666.3335
$ws.Range("N80").Value = -6662.3335

$ws.Range("H83").Value = 4499.75
$ws.Range("J83").Value = 4666.3335
$ws.Range("L83").Value = 23331.6675
$ws.Range("N83").Value = -33315.6675

$ws.Range("H132").Value = 84086.2
$ws.Range("I132").Value = 102258.336
$ws.Range("J132").Value = 11397.667
$ws.Range("K132").Value = 306775.008
$ws.Range("L132").Value = 34193.001
$ws.Range("M132").Value = -304245.008
$ws.Range("N132").Value = -39253.001


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 516
$ws.Range("I22").Value = 250
$ws.Range("K22").Value = 250
$ws.Range("M22").Value = 45

$ws.Range("H27").Value = 516
$ws.Range("I27").Value = 250
$ws.Range("K27").Value = 250
$ws.Range("M27").Value = -143

$ws.Range("H46").Value = 6633.5
$ws.Range("I46").Value = 1301
$ws.Range("J46").Value = 7700
$ws.Range("K46").Value = 1301
$ws.Range("L46").Value = 7700
$ws.Range("M46").Value = -1113
$ws.Range("N46").Value = -8076

$ws.Range("H55").Value = 1308.6923
$ws.Range("J55").Value = 3095
$ws.Range("L55").Value = 3095
$ws.Range("N55").Value = -3441

$ws.Range("H61").Value = 4628.2856
$ws.Range("I61").Value = 4049.5715
$ws.Range("J61").Value = 5785.7144
$ws.Range("K61").Value = 4049.5715
$ws.Range("L61").Value = 5785.7144
$ws.Range("M61").Value = -3847.5715
$ws.Range("N61").Value = -6189.7144

$ws.Range("H113").Value = 4628.2856
$ws.Range("I113").Value = 4049.5715
$ws.Range("J113").Value = 5785.7144
$ws.Range("K113").Value = 4049.5715
$ws.Range("L113").Value = 5785.7144
$ws.Range("M113").Value = -1879.5715
$ws.Range("N113").Value = -10125.7144

$ws.Range("H122").Value = 4833.3335
$ws.Range("I122").Value = 4000
$ws.Range("K122").Value = 12000
$ws.Range("M122").Value = -9550


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 54993.5
$ws.Range("J54").Value = 54993.5
$ws.Range("L54").Value = 54993.5
$ws.Range("N54").Value = -56033.5

$ws.Range("H81").Value = 5284.857
$ws.Range("I81").Value = 4623.75
$ws.Range("K81").Value = 9247.5
$ws.Range("M81").Value = -8186.5

$ws.Range("H84").Value = 5284.857
$ws.Range("I84").Value = 4623.75
$ws.Range("K84").Value = 46237.5
$ws.Range("M84").Value = -40933.5

$ws.Range("H113").Value = 535.619
$ws.Range("I113").Value = 490.13333
$ws.Range("K113").Value = 1470.39999
$ws.Range("M113").Value = 699.6000100000001

$ws.Range("H122").Value = 5999
$ws.Range("J122").Value = 5999
$ws.Range("L122").Value = 17997
$ws.Range("N122").Value = -22897

$ws.Range("H126").Value = 4338.0586
$ws.Range("I126").Value = 2931.5454
$ws.Range("J126").Value = 6916.6665
$ws.Range("K126").Value = 8794.636200000001
$ws.Range("L126").Value = 20749.9995
$ws.Range("M126").Value = -6324.636200000001
$ws.Range("N126").Value = -25689.9995

